$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2017")

# Update share counts (D column) for the affected holdings
$ws.Range("D2").Value = 27.971
$ws.Range("D3").Value = 14.606999999999999
$ws.Range("D6").Value = 25.283000000000001
$ws.Range("D11").Value = 15.205

# Record June dividends (M column) for the affected holdings
$ws.Range("M2").Value = 8.8800000000000008
$ws.Range("M3").Value = 6.95
$ws.Range("M6").Value = 3.02
$ws.Range("M11").Value = 5.7

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("J18").Select()
